$wb = $excel.ActiveWorkbook

foreach ($name in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 8378
    $ws.Range("F3").Value = 7929
    $ws.Range("F9").Value = 129
    $ws.Range("F11").Value = 233
    $ws.Range("F12").Value = 716
    $ws.Range("F13").Value = 135
    $ws.Range("F14").Value = 1908
    $ws.Range("F15").Value = 63
    $ws.Range("F16").Value = 58
    $ws.Range("F17").Value = 15
    $ws.Range("F20").Value = 12
}
